# Applies "Colocando header nos graficos" edit:
#  - Adds a header label in A1 (styled like the other header cells B1:E1) on every
#    data sheet, and de-bolds/un-borders the row-label cells A2:A12 (or A2:A3) that
#    previously shared that header style.
#  - Fixes a handful of missing-accent typos in the row labels.
#  - Removes the obsolete "Teto" row from the Emissoes sheet.
#  - Updates the Custo Total sheet: gives it a real header row and new values.

$wb = $excel.ActiveWorkbook

# xlPasteSpecial constant
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
#             "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# ---------------------------------------------------------------------------
$fixups = @{
    "A3"  = "Gás Natural"
    "A4"  = "Carvão"
    "A6"  = "Óleos Comb"
    "A8"  = "Eólica"
    "A11" = "Pot. Compl."
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # New header cell A1, styled like the neighbouring year headers.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)

    # Row labels A2:A12 lose the bold/border header styling ...
    for ($r = 2; $r -le 12; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.ClearFormats()
    }

    # ... and a few of them get their missing accents restored.
    foreach ($addr in $fixups.Keys) {
        $ws.Range($addr).Value = $fixups[$addr]
    }
}

# ---------------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial($xlPasteFormats)

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").ClearFormats()

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").ClearFormats()

# Row 4 ("Teto") is no longer used - remove it entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A1").Value = "Tipo Expansão"
$ws1.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)

# B1 used to read "Custo"; it now matches the other sheets' "2015" year header,
# keeping that header's text formatting (stored as text, not a number).
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4104)

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").ClearFormats()
$ws6.Range("B2").Value = 576

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").ClearFormats()
$ws6.Range("B3").Value = 99
